$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @(45754,45755,45756,45757,45758,45759,45760,45761,45762,45763,45764,45765,45766,45767,45768,45769,45770,45771,45772,45773,45774,45775,45776,45777,45778,45779,45780,45781,45782,45783)
$precs = @(0,0,0,0,11.3,4.5,3.6,4.9,0,2.8,0,12.7,4.8,0,0.9,0,0,0,0,0,0,0,0,5.3,22.7,13.6,0.9,2.8,19.4,0)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $precs[$i]
}
